$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "87.945.72"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "3.272.07"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'212.57"
$ws.Range("E5").Value = "  -3.05%  "
$ws.Range("D6").Value = "'627.59"
$ws.Range("E6").Value = "  -1.17%  "
$ws.Range("D7").Value = "'0.378"
$ws.Range("E7").Value = "  +15.17%  "
$ws.Range("D8").Value = "'0.715"
$ws.Range("E8").Value = "  +16.61%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").Value = "3.263.89"
$ws.Range("E10").Value = "  -1.43%  "
$ws.Range("E11").Value = "  -4.64%  "
$ws.Range("E12").Value = "  +11.84%  "
$ws.Range("D13").Value = "'0.0000263"
$ws.Range("E13").Value = "  -4.03%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "'34.49"
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("B15").Value = "Toncoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D15").Value = "'5.51"
$ws.Range("E15").Value = "  +1.76%  "
$ws.Range("D16").Value = "3.887.04"
$ws.Range("E16").Value = "  -1.33%  "
$ws.Range("D17").Value = "87.993.58"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").Value = "3.273.76"
$ws.Range("E18").Value = "  -1.64%  "
$ws.Range("D19").Value = "'3.18"
$ws.Range("E19").Value = "  -1.13%  "
$ws.Range("D20").Value = "'14.12"
$ws.Range("E20").Value = "  -2.92%  "
$ws.Range("D21").Value = "'438.67"
$ws.Range("E21").Value = "  -2.26%  "
$ws.Range("D22").Value = "'8.96"
$ws.Range("E22").Value = "  -1.30%  "
$ws.Range("D23").Value = "'5.34"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "'7.44"
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("D25").Value = "'12.34"
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("D26").Value = "'5.26"
$ws.Range("E26").Value = "  -1.90%  "
$ws.Range("D27").Value = "3.451.58"
$ws.Range("E27").Value = "  -1.65%  "
$ws.Range("D28").Value = "'77.29"
$ws.Range("E28").Value = "  -1.63%  "
$ws.Range("D29").Value = "'0.0000137"
$ws.Range("E29").Value = "  +6.14%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").Value = "'0.178"
$ws.Range("E31").Value = "  -16.83%  "
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  -0.18%  "
$ws.Range("D33").Value = "'568.04"
$ws.Range("E33").Value = "  -4.82%  "
$ws.Range("E34").Value = "  -5.15%  "
$ws.Range("D35").Value = "'1.38"
$ws.Range("E35").Value = "  -10.00%  "
$ws.Range("E36").Value = "  +8.56%  "
$ws.Range("E37").Value = "  -3.74%  "
$ws.Range("E38").Value = "  -7.33%  "
$ws.Range("D39").Value = "'22.86"
$ws.Range("E39").Value = "  -2.56%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("B41").Value = "WhiteBITCoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D41").Value = "'21.79"
$ws.Range("E41").Value = "  +1.81%  "
$ws.Range("E42").Value = "  +2.12%  "
$ws.Range("D43").Value = "'0.401"
$ws.Range("E43").Value = "  -3.84%  "
$ws.Range("E44").Value = "  -1.02%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Value = "'152.71"
$ws.Range("E46").Value = "  -2.86%  "
$ws.Range("D47").Value = "'0.136"
$ws.Range("E47").Value = "  +21.01%  "
$ws.Range("D48").Value = "'180.50"
$ws.Range("E48").Value = "  -4.48%  "
$ws.Range("D49").Value = "'44.77"
$ws.Range("E49").Value = "  -3.33%  "
$ws.Range("E50").Value = "  -3.09%  "
$ws.Range("E51").Value = "  -1.16%  "
